$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Single Robot")
$tbl = $ws.ListObjects.Item("Table1")

$row = $ws.Rows.Item(28)
$row.Insert()

$tbl.Resize($ws.Range("A1:F29"))

$ws.Cells.Item(28,1).Value = "zip ties"
$ws.Cells.Item(28,2).Value = 1

$url = "https://www.amazon.com/Nylon-Cable-Tie-Kit-Assorted/dp/B071SLNHZ3/ref=sr_1_8?crid=2NR8X8M4F90C0&dchild=1&keywords=zip+ties+assorted+sizes&qid=1598750768&s=hi&sprefix=zip+%2Ctools%2C236&sr=1-8"
$ws.Hyperlinks.Add($ws.Cells.Item(28,4), $url, "", "", $url) | Out-Null
$ws.Cells.Item(28,4).Style = "Hyperlink"

$ws.Cells.Item(28,5).Value = "Core"
$ws.Cells.Item(28,6).Value = "this is one large pack"

for ($r=1; $r -le 31; $r++) {
    $v = $ws.Cells.Item($r,1).Value2
    Write-Host "row$r : $v"
}
